# Wrote text introduction for CCPT section
# This script reproduces the authoring edit: it adds a new "SD_CPT / CPT Time"
# column (column I, rows 19-28) computed as H/G for each data row, and updates
# the saved window/selection state to match where the author had clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition the application window the way the author's Excel window was
#     positioned when the file was saved (best-effort; harmless if the host
#     does not persist on-screen window coordinates). ---
try {
    $win = $excel.Windows.Item(1)
    $win.Left   = -800
    $win.Top    = 4040
    $win.Width  = 38400
    $win.Height = 21160
} catch {
}

# --- Fill in the new "SD_CPT / CPT Time" ratio column (I19:I28) ---
# I19 = H19 / G19
$ws.Range("I19").Formula = "=H19/G19"
$ws.Range("I19").Style = "Normal"

# I20:I28 = H/G for each respective row (entered as one fill so the engine
# creates a shared formula, matching the authored workbook)
$ws.Range("I20:I28").Formula = "=H20/G20"
$ws.Range("I20:I28").Style = "Normal"

# --- Update the sheet selection to reflect the author's last click/selection ---
$ws.Range("I19:I28").Select()
